$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (row, column index, new value)
$updates = @(
    @(2, 7, 15.953202),
    @(2, 8, 47.859606),
    @(2, 9, 0.6210379196599995),
    @(2, 10, 0.6210379196599995),
    @(2, 13, 68.63737500000001),
    @(2, 14, 205.912125),
    @(2, 15, 0.5415701538216162),
    @(2, 16, 0.5415701538216162),
    @(2, 17, 1094.98590812475),
    @(2, 18, 9854.873173122751),
    @(2, 19, 0.3363356016793225),
    @(2, 20, 0.3363356016793225),
    @(3, 7, 15.953202),
    @(3, 8, 47.859606),
    @(3, 9, 0.6210379196599995),
    @(3, 10, 0.6210379196599995),
    @(3, 15, 0.08718851262838957),
    @(3, 16, 0.08718851262838957),
    @(3, 17, 176.284073272416),
    @(3, 18, 1586.556659451744),
    @(3, 19, 0.05414737250098466),
    @(3, 20, 0.05414737250098466),
    @(4, 7, 15.953202),
    @(4, 8, 47.859606),
    @(4, 9, 0.6210379196599995),
    @(4, 10, 0.6210379196599995),
    @(4, 13, 16.21089566666667),
    @(4, 14, 48.632687),
    @(4, 15, 0.1279089892319285),
    @(4, 16, 0.1279089892319285),
    @(4, 17, 258.615693171258),
    @(4, 18, 2327.541238541322),
    @(4, 19, 0.07943633257841017),
    @(4, 20, 0.07943633257841017),
    @(5, 7, 15.953202),
    @(5, 8, 47.859606),
    @(5, 9, 0.6210379196599995),
    @(5, 10, 0.6210379196599995),
    @(5, 13, 20.32546233333333),
    @(5, 14, 60.976387),
    @(5, 15, 0.1603741949973873),
    @(5, 16, 0.1603741949973873),
    @(5, 17, 324.256206347058),
    @(5, 18, 2918.305857123522),
    @(5, 19, 0.09959845642832454),
    @(5, 20, 0.09959845642832454),
    @(6, 7, 15.953202),
    @(6, 8, 47.859606),
    @(6, 9, 0.6210379196599995),
    @(6, 10, 0.6210379196599995),
    @(6, 13, 10.513928),
    @(6, 14, 31.541784),
    @(6, 15, 0.08295814932067838),
    @(6, 16, 0.08295814932067838),
    @(6, 17, 167.730817197456),
    @(6, 18, 1509.577354777104),
    @(6, 19, 0.05152015647295771),
    @(6, 20, 0.05152015647295771),
    @(7, 7, 0.7397413333333333),
    @(7, 9, 0.02879719185777549),
    @(7, 10, 0.02879719185777549),
    @(7, 13, 68.63737500000001),
    @(7, 14, 205.912125),
    @(7, 15, 0.5415701538216162),
    @(7, 16, 0.5415701538216162),
    @(7, 17, 50.773903299),
    @(7, 18, 456.965129691),
    @(7, 19, 0.01559569962404606),
    @(7, 20, 0.01559569962404606),
    @(8, 7, 0.7397413333333333),
    @(8, 9, 0.02879719185777549),
    @(8, 10, 0.02879719185777549),
    @(8, 15, 0.08718851262838957),
    @(8, 16, 0.08718851262838957),
    @(8, 17, 8.174196967352888),
    @(8, 18, 73.56777270617599),
    @(8, 19, 0.002510784325953815),
    @(8, 20, 0.002510784325953815),
    @(9, 7, 0.7397413333333333),
    @(9, 9, 0.02879719185777549),
    @(9, 10, 0.02879719185777549),
    @(9, 13, 16.21089566666667),
    @(9, 14, 48.632687),
    @(9, 15, 0.1279089892319285),
    @(9, 16, 0.1279089892319285),
    @(9, 17, 11.99186957498756),
    @(9, 18, 107.926826174888),
    @(9, 19, 0.003683419703245984),
    @(9, 20, 0.003683419703245984),
    @(10, 7, 0.7397413333333333),
    @(10, 9, 0.02879719185777549),
    @(10, 10, 0.02879719185777549),
    @(10, 13, 20.32546233333333),
    @(10, 14, 60.976387),
    @(10, 15, 0.1603741949973873),
    @(10, 16, 0.1603741949973873),
    @(10, 17, 15.03558460707644),
    @(10, 18, 135.320261463688),
    @(10, 19, 0.004618326462376061),
    @(10, 20, 0.004618326462376061),
    @(11, 7, 0.7397413333333333),
    @(11, 9, 0.02879719185777549),
    @(11, 10, 0.02879719185777549),
    @(11, 13, 10.513928),
    @(11, 14, 31.541784),
    @(11, 15, 0.08295814932067838),
    @(11, 16, 0.08295814932067838),
    @(11, 17, 7.777587117290666),
    @(11, 18, 69.99828405561598),
    @(11, 19, 0.002388961742153562),
    @(11, 20, 0.002388961742153562),
    @(12, 7, 5.607355000000001),
    @(12, 8, 16.822065),
    @(12, 9, 0.2182872180766656),
    @(12, 10, 0.2182872180766656),
    @(12, 13, 68.63737500000001),
    @(12, 14, 205.912125),
    @(12, 15, 0.5415701538216162),
    @(12, 16, 0.5415701538216162),
    @(12, 17, 384.8741278931251),
    @(12, 18, 3463.867151038125),
    @(12, 19, 0.1182178422710725),
    @(12, 20, 0.1182178422710725),
    @(13, 7, 5.607355000000001),
    @(13, 8, 16.822065),
    @(13, 9, 0.2182872180766656),
    @(13, 10, 0.2182872180766656),
    @(13, 15, 0.08718851262838957),
    @(13, 16, 0.08718851262838957),
    @(13, 17, 61.96169143250668),
    @(13, 18, 557.6552228925601),
    @(13, 19, 0.01903213786989339),
    @(13, 20, 0.01903213786989338),
    @(14, 7, 5.607355000000001),
    @(14, 8, 16.822065),
    @(14, 9, 0.2182872180766656),
    @(14, 10, 0.2182872180766656),
    @(14, 13, 16.21089566666667),
    @(14, 14, 48.632687),
    @(14, 15, 0.1279089892319285),
    @(14, 16, 0.1279089892319285),
    @(14, 17, 90.90024687096169),
    @(14, 18, 818.1022218386552),
    @(14, 19, 0.02792089742643585),
    @(14, 20, 0.02792089742643585),
    @(15, 7, 5.607355000000001),
    @(15, 8, 16.822065),
    @(15, 9, 0.2182872180766656),
    @(15, 10, 0.2182872180766656),
    @(15, 13, 20.32546233333333),
    @(15, 14, 60.976387),
    @(15, 15, 0.1603741949973873),
    @(15, 16, 0.1603741949973873),
    @(15, 17, 113.9720828421284),
    @(15, 18, 1025.748745579155),
    @(15, 19, 0.03500763687726438),
    @(15, 20, 0.03500763687726438),
    @(16, 7, 5.607355000000001),
    @(16, 8, 16.822065),
    @(16, 9, 0.2182872180766656),
    @(16, 10, 0.2182872180766656),
    @(16, 13, 10.513928),
    @(16, 14, 31.541784),
    @(16, 15, 0.08295814932067838),
    @(16, 16, 0.08295814932067838),
    @(16, 17, 58.95532674044001),
    @(16, 18, 530.5979406639601),
    @(16, 19, 0.01810870363199951),
    @(16, 20, 0.01810870363199951),
    @(17, 7, 1.453021),
    @(17, 8, 4.359063),
    @(17, 9, 0.05656426459479998),
    @(17, 10, 0.05656426459479998),
    @(17, 13, 68.63737500000001),
    @(17, 14, 205.912125),
    @(17, 15, 0.5415701538216162),
    @(17, 16, 0.5415701538216162),
    @(17, 17, 99.731547259875),
    @(17, 18, 897.583925338875),
    @(17, 19, 0.03063351747741242),
    @(17, 20, 0.03063351747741243),
    @(18, 7, 1.453021),
    @(18, 8, 4.359063),
    @(18, 9, 0.05656426459479998),
    @(18, 10, 0.05656426459479998),
    @(18, 15, 0.08718851262838957),
    @(18, 16, 0.08718851262838957),
    @(18, 17, 16.05599054223467),
    @(18, 18, 144.503914880112),
    @(18, 19, 0.004931754097939287),
    @(18, 20, 0.004931754097939288),
    @(19, 7, 1.453021),
    @(19, 8, 4.359063),
    @(19, 9, 0.05656426459479998),
    @(19, 10, 0.05656426459479998),
    @(19, 13, 16.21089566666667),
    @(19, 14, 48.632687),
    @(19, 15, 0.1279089892319285),
    @(19, 16, 0.1279089892319285),
    @(19, 17, 23.55477183247567),
    @(19, 18, 211.992946492281),
    @(19, 19, 0.007235077910968225),
    @(19, 20, 0.007235077910968226),
    @(20, 7, 1.453021),
    @(20, 8, 4.359063),
    @(20, 9, 0.05656426459479998),
    @(20, 10, 0.05656426459479998),
    @(20, 13, 20.32546233333333),
    @(20, 14, 60.976387),
    @(20, 15, 0.1603741949973873),
    @(20, 16, 0.1603741949973873),
    @(20, 17, 29.53332360504233),
    @(20, 18, 265.799912445381),
    @(20, 19, 0.009071448400010264),
    @(20, 20, 0.009071448400010266),
    @(21, 7, 1.453021),
    @(21, 8, 4.359063),
    @(21, 9, 0.05656426459479998),
    @(21, 10, 0.05656426459479998),
    @(21, 13, 10.513928),
    @(21, 14, 31.541784),
    @(21, 15, 0.08295814932067838),
    @(21, 16, 0.08295814932067838),
    @(21, 17, 15.276958176488),
    @(21, 18, 137.492623588392),
    @(21, 19, 0.004692466708469778),
    @(21, 20, 0.004692466708469779),
    @(22, 7, 1.934648333333333),
    @(22, 8, 5.803945),
    @(22, 9, 0.07531340581075942),
    @(22, 10, 0.07531340581075942),
    @(22, 13, 68.63737500000001),
    @(22, 14, 205.912125),
    @(22, 15, 0.5415701538216162),
    @(22, 16, 0.5415701538216162),
    @(22, 17, 132.789183148125),
    @(22, 18, 1195.102648333125),
    @(22, 19, 0.04078749276976278),
    @(22, 20, 0.04078749276976278),
    @(23, 7, 1.934648333333333),
    @(23, 8, 5.803945),
    @(23, 9, 0.07531340581075942),
    @(23, 10, 0.07531340581075942),
    @(23, 15, 0.08718851262838957),
    @(23, 16, 0.08718851262838957),
    @(23, 17, 21.37800853707556),
    @(23, 18, 192.40207683368),
    @(23, 19, 0.006566463833618426),
    @(23, 20, 0.006566463833618426),
    @(24, 7, 1.934648333333333),
    @(24, 8, 5.803945),
    @(24, 9, 0.07531340581075942),
    @(24, 10, 0.07531340581075942),
    @(24, 13, 16.21089566666667),
    @(24, 14, 48.632687),
    @(24, 15, 0.1279089892319285),
    @(24, 16, 0.1279089892319285),
    @(24, 17, 31.36238228335723),
    @(24, 18, 282.261440550215),
    @(24, 19, 0.009633261612868288),
    @(24, 20, 0.009633261612868288),
    @(25, 7, 1.934648333333333),
    @(25, 8, 5.803945),
    @(25, 9, 0.07531340581075942),
    @(25, 10, 0.07531340581075942),
    @(25, 13, 20.32546233333333),
    @(25, 14, 60.976387),
    @(25, 15, 0.1603741949973873),
    @(25, 16, 0.1603741949973873),
    @(25, 17, 39.32262182741278),
    @(25, 18, 353.903596446715),
    @(25, 19, 0.01207832682941209),
    @(25, 20, 0.01207832682941209),
    @(26, 7, 1.934648333333333),
    @(26, 8, 5.803945),
    @(26, 9, 0.07531340581075942),
    @(26, 10, 0.07531340581075942),
    @(26, 13, 10.513928),
    @(26, 14, 31.541784),
    @(26, 15, 0.08295814932067838),
    @(26, 16, 0.08295814932067838),
    @(26, 17, 20.34075328198666),
    @(26, 18, 183.06677953788),
    @(26, 19, 0.006247860765097826),
    @(26, 20, 0.006247860765097826)
)

foreach ($u in $updates) {
    $ws.Cells.Item([int]$u[0], [int]$u[1]).Value = [double]$u[2]
}
